$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 32: "Guide needs up/down/left/right buttons..." is now Done.
# Remove the old comment in C32 and mark B32 = Done (matching formatting
# of the other Status cells: vertical-top + wrap text).
$ws.Range("C32").Clear()
$ws.Range("B32").Value = "Done"
$ws.Range("B32").WrapText = $true
$ws.Range("B32").VerticalAlignment = -4160

# Row 33: "Guide rows need to be a little bigger, also guide font" is now Done.
# Remove the old comment in C33 and mark B33 = Done (same formatting).
$ws.Range("C33").Clear()
$ws.Range("B33").Value = "Done"
$ws.Range("B33").WrapText = $true
$ws.Range("B33").VerticalAlignment = -4160

# Update selection to reflect where the author finished editing.
$ws.Range("C33").Select()
